$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the duplicate-suffixed headers back to their plain forms.
# Apply in an order that frees up the old shared-string slots before new
# ones are appended, so the resulting shared string table lands in the
# same order Excel produced: Factor [time], Term Source REF (PATO:0000165),
# Term Accession Number (PATO:0000165), Parameter [temperature unit],
# Term Source REF (UO:0000005), Term Accession Number (UO:0000005).
$ws.Range("F2").Value = "Factor [time]"
$ws.Range("H2").Value = "Term Source REF (PATO:0000165)"
$ws.Range("I2").Value = "Term Accession Number (PATO:0000165)"
$ws.Range("J2").Value = "Parameter [temperature unit]"
$ws.Range("L2").Value = "Term Source REF (UO:0000005)"
$ws.Range("M2").Value = "Term Accession Number (UO:0000005)"

# Update the view: move the selection (this also resets the scrolled
# topLeftCell back to the default).
$ws.Range("J2").Select()
